$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (number of reps/columns) change
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values change, C2 is cleared entirely (cell removed)
$ws.Range("B2").Value = 0.89786533337632157
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.4391211244529551
$ws.Range("E2").Value = 0.99779364836851714

# Row 3 data values change, B3 cleared entirely, D3 newly populated
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 0.55039406598732965
$ws.Range("D3").Value = 1.9446659266568205
$ws.Range("E3").Value = 0.20276339621892436

# Selection narrows to B1:E3 as shown in sheetView
$ws.Range("B1:E3").Select()
